# Apply updates to the "Inscricoes" sheet (Resumo de Inscricoes) as described
# in the commit diff: a handful of registration counts were incremented.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Column E = "Inscritos"
$ws.Range("E9").Value  = 19
$ws.Range("E12").Value = 2
$ws.Range("E18").Value = 87
$ws.Range("E32").Value = 15
$ws.Range("E36").Value = 72
$ws.Range("E47").Value = 47
$ws.Range("E49").Value = 53
$ws.Range("E63").Value = 20

# Row 70: Inscritos, Pagos, Inscricoes homologadas all bump by 1
$ws.Range("E70").Value = 29
$ws.Range("F70").Value = 13
$ws.Range("H70").Value = 13

# Row 77: Pagos and Inscricoes homologadas bump by 1
$ws.Range("F77").Value = 15
$ws.Range("H77").Value = 15
